$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 434 (shifts existing rows 434..531 down to 435..532,
# and bumps the sheet dimension from A1:R531 to A1:R532).
$ws.Rows(434).Insert()

# Populate the newly inserted row 434 with a new price observation
# (same market/region/variety context as the old row 434, now at 435,
# but a newer date and "1a (cosecha)" quality).
$ws.Cells.Item(434, 1).Value  = 8
$ws.Cells.Item(434, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(434, 3).Value  = "Coquimbo"
$ws.Cells.Item(434, 4).Value  = 44889
$ws.Cells.Item(434, 5).Value  = 4
$ws.Cells.Item(434, 6).Value  = 100114001
$ws.Cells.Item(434, 7).Value  = "Papa"
$ws.Cells.Item(434, 8).Value  = "Cardinal"
$ws.Cells.Item(434, 9).Value  = "1a (cosecha)"
$ws.Cells.Item(434, 10).Value = 2000
$ws.Cells.Item(434, 11).Value = 11500
$ws.Cells.Item(434, 12).Value = 12000
$ws.Cells.Item(434, 13).Value = 11750
$ws.Cells.Item(434, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(434, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(434, 16).Value = 470
$ws.Cells.Item(434, 17).Value = 25
$ws.Cells.Item(434, 18).Value = "Hortaliza"
